# SCRUM_TEAM_D.xlsx edit script
# - Switches the active tab from "Backlog" to "Sprint-Backlog"
# - Replaces the short "Umsetzung ..." placeholder descriptions with the
#   detailed versions (on both the Backlog and Sprint-Backlog sheets) and
#   reorders two Backlog rows
# - Fills in the first sprint's rows 5-7 on the Sprint-Backlog sheet with
#   the corresponding Backlog items plus estimation/tracking columns
#   (D/E/F) and a new "Focus-Faktor" note in D1

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Backlog
$ws2 = $wb.Worksheets.Item(2)   # Sprint-Backlog

# ---------------------------------------------------------------
# Backlog sheet (ws1): update descriptions to the longer wording,
# and swap the "Spezifikation eines weiteren Spiels" /
# "Erzeugung der TicTacToe-Spieler in einer Fabrik." rows.
# ---------------------------------------------------------------
$ws1.Range("B1").Value = "Team D"

$ws1.Range("C6").Value = "Umsetzung des Painters im Framework: Standard TicTacToe. Kreise und Kreuz und Bunt. Mit Ausblick auf Timeranzeige. Mit ausblick auf Gewinneranzeige im Spielfeld."

$ws1.Range("B7").Value = "Implementierung der TicTacToe-Regeln"
$ws1.Range("C7").Value = "Umsetzung der Rules im Framework: TicTactoe mit 30s Limit pro Zug. Spielerwechsel falls nicht eingehalten. Sonst Standard TicTacToe"

$ws1.Range("B8").Value = "Implementierung eines menschlichen TicTacToe-Spielers"
$ws1.Range("C8").Value = "Umsetzung eines menschlichen Players im Framework: Standard TicTacToe mit Mauseingabe. Ausblick auf Farbenwahl durch Spieler. Mit Ausblick auf zufällig beginnenden Spieler."

$ws1.Range("B9").Value = "Implementierung eines PC gesteurten TicTacToe-Spielers"
$ws1.Range("C9").Value = "Umsetzung eines PC gesteuerten Players im Framework."

$ws1.Range("B10").Value = "Erzeugung der TicTacToe-Spieler in einer Fabrik."
$ws1.Range("C10").Value = "Umsetung einer Fabrik zur Erzeugung der menschlichen und PC gesteuerten TicTacToe-Spieler."

$ws1.Range("B11").Value = "Spezifikation eines weiteren Spiels"
$ws1.Range("C11").Value = "Im Entwicklungsteam kann beschlossen werden welches weitere Spiel umgesetzt werden soll, bzw. kann."

# Backlog sheet selection: full row 7 selected, tab not active
$ws1.Rows.Item(7).Select()

# ---------------------------------------------------------------
# Sprint-Backlog sheet (ws2): new "Focus-Faktor" note, plan/track
# columns for the first two rows, and three freshly planned items.
# ---------------------------------------------------------------
$ws2.Range("D1").Value = "Focus-Faktor: 0,5"

# Build the centered (no-wrap) style used by column E from scratch on E3,
# then propagate it via copy/paste-format so only one new style gets
# created in xl/styles.xml.
$e3 = $ws2.Range("E3")
$e3.VerticalAlignment = -4108   # xlCenter
$e3.HorizontalAlignment = -4108 # xlCenter
$e3.Copy()
$ws2.Range("E4:E5").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws2.Range("D3").Value = "Due: 04.11.2021"
$ws2.Range("E3").Value = "Done"

$ws2.Range("D4").Value = "Due: 04.11.2021"

$ws2.Range("A5").Value = 2
$ws2.Range("B5").Value = "Spezifikation einer Umsetzungsidee für das Spiel TicTacToe"
$ws2.Range("C5").Value = "Spezifikation: Wie soll das TicTacToe aussehen? Wie soll das Aussehen technisch erreicht werden? Wie sollen Regeln umgesetzt werden? Soll es Sonderregeln geben, bzw. Sonderspielfelder, etc.? Welche Spieler (PC, ...) soll es geben? Etc.?"
$ws2.Range("D5").Value = "10min"
$ws2.Range("E5").Value = "Done"

$ws2.Range("A6").Value = 2
$ws2.Range("B6").Value = "Implementierung der TicTacToe-Darstellung"
$ws2.Range("C6").Value = "Umsetzung des Painters im Framework: Standard TicTacToe. Kreise und Kreuz und Bunt. Mit Ausblick auf Timeranzeige. Mit ausblick auf Gewinneranzeige im Spielfeld."
$ws2.Range("D6").Value = "180min"

$ws2.Range("A7").Value = 2
$ws2.Range("B7").Value = "Implementierung der TicTacToe-Regeln"
$ws2.Range("C7").Value = "Umsetzung der Rules im Framework: TicTactoe mit 30s Limit pro Zug. Spielerwechsel falls nicht eingehalten. Sonst Standard TicTacToe"
$ws2.Range("D7").Value = "180min"
$ws2.Range("F7").Value = "Überplant um 90mins"

# Sprint-Backlog becomes the active sheet/tab, selection on E4
$ws2.Activate()
$ws2.Range("E4").Select()
